# Natmi following Dr Hou advice
# Replace the single "sCs -> Efnb3/Rhbdl2 -> ECs" LR-pair result row with the
# recomputed 4-cluster result set (M1, M2, Neutro, sCs), each recalculated
# against the larger background of sending clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("M1",     "Efnb3", "Rhbdl2", "ECs", 1, 0.3333333333333333, 0.06315533333333333, 0.189466,            0.08692959137136654, 0.08692959137136655, 3, 1, 1.021610666666667, 3.064832, 1, 1, 0.06452016219022222, 0.580681459712,     0.08692959137136654, 0.08692959137136655),
    @("M2",     "Efnb3", "Rhbdl2", "ECs", 1, 0.3333333333333333, 0.028136,            0.084408,            0.03872754451180849, 0.03872754451180849, 3, 1, 1.021610666666667, 3.064832, 1, 1, 0.02874403771733333, 0.258696339456,     0.03872754451180849, 0.03872754451180849),
    @("Neutro", "Efnb3", "Rhbdl2", "ECs", 1, 0.3333333333333333, 0.073382,            0.220146,            0.1010059948594516,  0.1010059948594516,  3, 1, 1.021610666666667, 3.064832, 1, 1, 0.07496783394133334, 0.6747105054720001, 0.1010059948594516,  0.1010059948594516),
    @("sCs",    "Efnb3", "Rhbdl2", "ECs", 3, 1,                   0.5618379999999999, 1.685514,            0.7733368692573733,  0.7733368692573734,  3, 1, 1.021610666666667, 3.064832, 1, 1, 0.5739796937386666,  5.165817243648,     0.7733368692573733,  0.7733368692573734)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowData = $rows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
